$p = $ppt.ActivePresentation

# --- Slide 10: fix typo "Sate" -> "State" in the title ---
# (the title shape has a trailing empty paragraph that must be preserved,
#  so only the first paragraph's characters are touched)
$s10 = $p.Slides.Item(10)
$titleShape10 = $s10.Shapes.Item(1)
$tr10 = $titleShape10.TextFrame.TextRange
$origLen10 = $tr10.Paragraphs(1, 1).Characters().Count
$tr10.Characters(1, $origLen10).Text = "State with most alcoholic (ABV) beer and most bitter beer"

# --- Slide 12: add a comma after "future research" ---
$s12 = $p.Slides.Item(12)
$bodyShape12 = $s12.Shapes.Item(4)
$bodyShape12.TextFrame.TextRange.Text = "There is evidence to suggest a positive correlation between IBU and ABV. Using this data would help the company create custom beer flavors for its customers. However, because this is an observational study for future research, we will run an experiment to see if there is an actual causal relationship."

# --- Slide 13: lower-case "findings"/"business" in the title, and resize/reposition the SmartArt diagram ---
$s13 = $p.Slides.Item(13)
$titleShape13 = $s13.Shapes.Item(1)
$titleShape13.TextFrame.TextRange.Text = "Overall findings and business suggestions"

$diagramShape = $s13.Shapes.Item(2)
$diagramShape.Left = 42.359762191772745
$diagramShape.Top = 124.61464309692427
$diagramShape.Width = 883.7319641113286
$diagramShape.Height = 384.0092926025393
